# Fills in the first empty "category / amount" row that follows the
# already-populated "Chi phí kính doanh, quảng cáo, tiếp thị" / "5000000"
# row in the cost-estimation table (3rd table in the document):
#   Col 1 (category, italic): "Chi phí kính doanh, quảng cáo, tiếp thị"
#   Col 2 (amount):           "5000000"

$d = $word.ActiveDocument

$table = $d.Tables.Item(3)
$row = $table.Rows.Item(6)

# --- Column 1: category text, italic ---
$cell1 = $row.Cells.Item(1)
$cell1Range = $cell1.Range
$cell1Range.Text = "Chi phí kính doanh, quảng cáo, tiếp thị"

# Re-grab the cell range (it now contains the new text plus the trailing
# end-of-cell mark) and italicize only the text, not the paragraph mark,
# so no extra paragraph-mark run-properties get introduced.
$textRange = $cell1.Range
$textOnly = $d.Range($textRange.Start, $textRange.End - 1)
$textOnly.Font.Italic = $true

# --- Column 2: amount ---
$cell2 = $row.Cells.Item(2)
$cell2Range = $cell2.Range
$cell2Range.Text = "5000000"
